$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range('D2').Value = '64.076.13'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '3.305.89'
$ws.Range('E3').Value = '  +5.87%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextCell 'D5' '599.50'
$ws.Range('E5').Value = '  +1.19%  '
Set-TextCell 'D6' '143.50'
$ws.Range('E6').Value = '  +5.12%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.300.71'
$ws.Range('E8').Value = '  +6.01%  '
Set-TextCell 'D9' '0.523'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  +2.56%  '
Set-TextCell 'D11' '5.47'
$ws.Range('E11').Value = '  +3.43%  '
Set-TextCell 'D12' '0.473'
$ws.Range('E12').Value = '  +3.01%  '
$ws.Range('E13').Value = '  -0.04%  '
Set-TextCell 'D14' '34.97'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').Value = '3.847.16'
$ws.Range('E15').Value = '  +5.86%  '
$ws.Range('E16').Value = '  +1.30%  '
$ws.Range('D17').Value = '3.302.62'
$ws.Range('E17').Value = '  +5.51%  '
$ws.Range('D18').Value = '64.139.30'
$ws.Range('E18').Value = '  +1.76%  '
Set-TextCell 'D19' '6.90'
$ws.Range('E19').Value = '  +2.67%  '
Set-TextCell 'D20' '483.78'
$ws.Range('E20').Value = '  +1.45%  '
Set-TextCell 'D21' '14.31'
$ws.Range('E21').Value = '  +1.17%  '
Set-TextCell 'D22' '0.744'
$ws.Range('E22').Value = '  +6.41%  '
Set-TextCell 'D23' '8.04'
$ws.Range('E23').Value = '  +4.63%  '
Set-TextCell 'D24' '13.56'
$ws.Range('E24').Value = '  +3.81%  '
Set-TextCell 'D25' '84.53'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('E26').Value = '  -0.09%  '
Set-TextCell 'D27' '2.78'
$ws.Range('E27').Value = '  +2.68%  '
Set-TextCell 'D28' '7.35'
$ws.Range('E28').Value = '  +2.45%  '
Set-TextCell 'D29' '8.31'
$ws.Range('E29').Value = '  +3.87%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +3.21%  '
Set-TextCell 'D32' '28.76'
$ws.Range('E32').Value = '  +6.08%  '
$ws.Range('E33').Value = '  -1.62%  '
Set-TextCell 'D34' '2.57'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('E35').Value = '  +1.90%  '
$ws.Range('E36').Value = '  +2.92%  '
Set-TextCell 'D37' '53.35'
$ws.Range('E37').Value = '  +2.55%  '
$ws.Range('D38').Value = '0.0₃0739'
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('E39').Value = '  +3.22%  '
Set-TextCell 'D40' '430.55'
$ws.Range('E40').Value = '  +2.51%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D41' '8.47'
$ws.Range('E41').Value = '  +2.38%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.011.53'
$ws.Range('E42').Value = '  +4.48%  '
$ws.Range('E43').Value = '  +3.40%  '
$ws.Range('E44').Value = '  -4.98%  '
Set-TextCell 'D45' '0.271'
$ws.Range('E45').Value = '  +2.88%  '
$ws.Range('E46').Value = '  +5.10%  '
Set-TextCell 'D47' '26.36'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell 'D49' '35.74'
$ws.Range('E49').Value = '  +15.32%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell 'D50' '2.33'
$ws.Range('E50').Value = '  +2.61%  '
$ws.Range('E51').Value = '  +1.54%  '
